$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")

# Copy formatting (style) from the last existing data row (72)
# down into the two new rows (73, 74) before writing their content.
$ws.Range("A72:C72").Copy($ws.Range("A73:C74"))
$ws.Rows.Item(73).RowHeight = 20
$ws.Rows.Item(74).RowHeight = 20

# Formula for the new rows (shares the ROW()-2 pattern used throughout column A)
$ws.Range("A73:A74").Formula = "=ROW()-2"

# Row 73: "命中率が<val1>%上昇した" / blue
$ws.Range("B73").Value = "命中率が<val1>%上昇した"
$ws.Range("C73").Value = "blue"

# Row 74: "回避率が<val1>%上昇した" / blue
$ws.Range("B74").Value = "回避率が<val1>%上昇した"
$ws.Range("C74").Value = "blue"
